# Updated with Excel Sample
# Adds a new "Captured_Values" worksheet after "Contact_Information", fills in
# sample captured data, sizes its columns, and leaves selection/active-sheet
# state matching the authored workbook (Captured_Values active, with
# Contact_Information's selection parked on C29).

$wb = $excel.ActiveWorkbook

$wsContacts = $wb.Worksheets.Item("Contact_Information")

# Insert the new sheet right after Contact_Information.
$wsCaptured = $wb.Worksheets.Add([Type]::Missing, $wsContacts)
$wsCaptured.Name = "Captured_Values"

# Header row
$wsCaptured.Range("A1").Value = "Number"
$wsCaptured.Range("B1").Value = "Text1"

# Data rows
$wsCaptured.Range("A2").Value = 123456789
$wsCaptured.Range("B2").Value = "Real Programmers Count 0123456789 From Zero"

$wsCaptured.Range("A3").Value = 123456789
$wsCaptured.Range("B3").Value = "Real Programmers Count 0123456789 From Zero"

$wsCaptured.Range("A4").Value = 123456789
$wsCaptured.Range("B4").Value = "Real Programmers Count 0123456789 From Zero"

$wsCaptured.Range("A5").Value = 123456789
$wsCaptured.Range("B5").Value = "Real Programmers Count 0123456789 From Zero"

# Column widths to roughly match the authored widths (~24.57 / ~61.57 chars).
$wsCaptured.Columns.Item(1).ColumnWidth = 23.6666666666666667
$wsCaptured.Columns.Item(2).ColumnWidth = 60.6666666666666667

# Restore selection on the original sheet, then make the new sheet active
# with its own selection, matching the saved workbook view state.
$null = $wsContacts.Range("C29").Select()
$null = $wsCaptured.Select()
$null = $wsCaptured.Range("A2:XFD8").Select()
